# Weekly update: insert 3 new rows (638-640) for date 44578 (2022-01-17)
# and shift all existing "Betarraga" rows down by 3. Net effect: the used
# range grows from A1:R693 to A1:R696.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows above the existing data block (old row 638).
$ws.Range("A638:A640").EntireRow.Insert()

# Final state for every row in the affected block (638-696): row number,
# Fecha (date serial in col D), Clasificacion (col I), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M, P).
$data = @(
    @(638, 44578, "Primera", 54000, 80, 90, 85, 85),
    @(639, 44578, "Segunda", 41000, 65, 70, 67, 67),
    @(640, 44578, "Tercera", 13000, 50, 50, 50, 50),
    @(641, 44490, "Primera", 43000, 90, 100, 94, 94),
    @(642, 44490, "Segunda", 36000, 70, 80, 75, 75),
    @(643, 44490, "Tercera", 13000, 60, 60, 60, 60),
    @(644, 44427, "Primera", 42000, 110, 120, 114, 114),
    @(645, 44427, "Segunda", 37000, 85, 90, 87, 87),
    @(646, 44427, "Tercera", 11000, 65, 65, 65, 65),
    @(647, 44491, "Primera", 46000, 90, 100, 95, 95),
    @(648, 44491, "Segunda", 34000, 80, 85, 82, 82),
    @(649, 44491, "Tercera", 12000, 60, 60, 60, 60),
    @(650, 44293, "Primera", 43000, 120, 150, 136, 136),
    @(651, 44293, "Segunda", 20000, 100, 120, 111, 111),
    @(652, 44266, "Primera", 24000, 120, 120, 120, 120),
    @(653, 44266, "Segunda", 13000, 95, 95, 95, 95),
    @(654, 44533, "Primera", 55000, 70, 75, 72, 72),
    @(655, 44533, "Segunda", 43000, 50, 60, 54, 54),
    @(656, 44533, "Tercera", 13000, 35, 35, 35, 35),
    @(657, 44264, "Primera", 27000, 120, 120, 120, 120),
    @(658, 44264, "Segunda", 15000, 90, 90, 90, 90),
    @(659, 44494, "Primera", 39000, 100, 110, 104, 104),
    @(660, 44494, "Tercera", 9000, 60, 60, 60, 60),
    @(661, 44571, "Primera", 56000, 75, 80, 77, 77),
    @(662, 44571, "Segunda", 46000, 60, 65, 62, 62),
    @(663, 44571, "Tercera", 15000, 50, 50, 50, 50),
    @(664, 44390, "Primera", 38000, 110, 120, 114, 114),
    @(665, 44390, "Segunda", 31000, 90, 95, 92, 92),
    @(666, 44390, "Tercera", 9000, 70, 70, 70, 70),
    @(667, 44279, "Primera", 59000, 110, 130, 120, 120),
    @(668, 44279, "Segunda", 27000, 100, 100, 100, 100),
    @(669, 44481, "Primera", 43000, 90, 100, 94, 94),
    @(670, 44481, "Segunda", 38000, 80, 85, 82, 82),
    @(671, 44481, "Tercera", 12000, 60, 60, 60, 60),
    @(672, 44277, "Primera", 37000, 100, 130, 119, 119),
    @(673, 44277, "Segunda", 14000, 90, 100, 96, 96),
    @(674, 44525, "Primera", 63000, 80, 85, 82, 82),
    @(675, 44525, "Segunda", 52000, 60, 65, 62, 62),
    @(676, 44525, "Tercera", 18000, 40, 40, 40, 40),
    @(677, 44327, "Primera", 33000, 100, 100, 100, 100),
    @(678, 44327, "Segunda", 23000, 75, 75, 75, 75),
    @(679, 44327, "Tercera", 12000, 60, 60, 60, 60),
    @(680, 44354, "Primera", 43000, 100, 110, 104, 104),
    @(681, 44354, "Segunda", 30000, 75, 80, 77, 77),
    @(682, 44354, "Tercera", 9000, 50, 50, 50, 50),
    @(683, 44503, "Primera", 45000, 90, 100, 95, 95),
    @(684, 44503, "Segunda", 39000, 75, 80, 77, 77),
    @(685, 44503, "Tercera", 13000, 60, 60, 60, 60),
    @(686, 44384, "Primera", 59000, 100, 120, 108, 108),
    @(687, 44384, "Segunda", 39000, 70, 90, 78, 78),
    @(688, 44384, "Tercera", 8000, 50, 50, 50, 50),
    @(689, 44512, "Primera", 54000, 80, 90, 84, 84),
    @(690, 44512, "Segunda", 44000, 70, 75, 72, 72),
    @(691, 44512, "Tercera", 15000, 50, 50, 50, 50),
    @(692, 44312, "Primera", 24000, 100, 100, 100, 100),
    @(693, 44312, "Segunda", 1000, 80, 80, 80, 80),
    @(694, 44511, "Primera", 50000, 90, 100, 94, 94),
    @(695, 44511, "Segunda", 43000, 75, 85, 79, 79),
    @(696, 44511, "Tercera", 17000, 60, 60, 60, 60)
)

foreach ($row in $data) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 100114014
    $ws.Cells.Item($r, 7).Value = "Betarraga"
    $ws.Cells.Item($r, 8).Value = "Sin especificar"
    $ws.Cells.Item($r, 9).Value = $row[2]
    $ws.Cells.Item($r, 10).Value = $row[3]
    $ws.Cells.Item($r, 11).Value = $row[4]
    $ws.Cells.Item($r, 12).Value = $row[5]
    $ws.Cells.Item($r, 13).Value = $row[6]
    $ws.Cells.Item($r, 14).Value = "$/unidad"
    $ws.Cells.Item($r, 15).Value = "Región Metropolitana"
    $ws.Cells.Item($r, 16).Value = $row[7]
    $ws.Cells.Item($r, 17).Value = 1
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}
